$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Schedule")
Write-Host ($ws.Name)
